$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '29.592.65'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").Value = '1.919.98'
$ws.Range("E3").Value = '  +0.21%  '

Set-TextValue $ws.Range("D4") '1.010'
$ws.Range("E4").Value = '  +0.19%  '

Set-TextValue $ws.Range("D5") '339.37'
$ws.Range("E5").Value = '  +4.39%  '

Set-TextValue $ws.Range("D6") '1.010'
$ws.Range("E6").Value = '  +0.36%  '

Set-TextValue $ws.Range("D7") '0.4806'
$ws.Range("E7").Value = '  -0.28%  '

Set-TextValue $ws.Range("D8") '0.4048'
$ws.Range("E8").Value = '  -0.69%  '

Set-TextValue $ws.Range("D9") '0.08086'
$ws.Range("E9").Value = '  -1.79%  '

Set-TextValue $ws.Range("D10") '0.9991'
$ws.Range("E10").Value = '  -1.62%  '

Set-TextValue $ws.Range("D11") '23.41'
$ws.Range("E11").Value = '  -0.27%  '

$ws.Range("D12").Value = '1.933.51'
$ws.Range("E12").Value = '  +0.90%  '

Set-TextValue $ws.Range("D13") '5.990'
$ws.Range("E13").Value = '  -1.66%  '

Set-TextValue $ws.Range("D14") '7.188'
$ws.Range("E14").Value = '  -0.75%  '

Set-TextValue $ws.Range("D15") '89.82'
$ws.Range("E15").Value = '  -1.62%  '

Set-TextValue $ws.Range("D16") '0.06846'
$ws.Range("E16").Value = '  +0.48%  '

Set-TextValue $ws.Range("D17") '1.011'
$ws.Range("E17").Value = '  +0.28%  '

Set-TextValue $ws.Range("D18") '0.00001027'
$ws.Range("E18").Value = '  -1.21%  '

Set-TextValue $ws.Range("D19") '17.52'
$ws.Range("E19").Value = '  -0.91%  '

Set-TextValue $ws.Range("D20") '1.009'
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("D21").Value = '29.621.69'
$ws.Range("E21").Value = '  +0.60%  '

Set-TextValue $ws.Range("D22") '5.542'
$ws.Range("E22").Value = '  -2.04%  '

Set-TextValue $ws.Range("D23") '11.74'
$ws.Range("E23").Value = '  -0.59%  '

Set-TextValue $ws.Range("D24") '2.160'
$ws.Range("E24").Value = '  -0.70%  '

$ws.Range("D25").Value = '2.178.48'
$ws.Range("E25").Value = '  +1.16%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D26") '156.78'
$ws.Range("E26").Value = '  +0.55%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D27") '6.552'
$ws.Range("E27").Value = '  -0.55%  '

Set-TextValue $ws.Range("D28") '19.82'
$ws.Range("E28").Value = '  -1.08%  '

Set-TextValue $ws.Range("D29") '2.066'
$ws.Range("E29").Value = '  -2.29%  '

Set-TextValue $ws.Range("D30") '120.12'
$ws.Range("E30").Value = '  -0.20%  '

Set-TextValue $ws.Range("D31") '0.9981'
$ws.Range("E31").Value = '  -2.31%  '

Set-TextValue $ws.Range("D32") '0.09569'
$ws.Range("E32").Value = '  -0.15%  '

Set-TextValue $ws.Range("D33") '5.518'
$ws.Range("E33").Value = '  -2.94%  '

Set-TextValue $ws.Range("D34") '1.396'
$ws.Range("E34").Value = '  +1.67%  '

Set-TextValue $ws.Range("D35") '3.538'
$ws.Range("E35").Value = '  -0.32%  '

Set-TextValue $ws.Range("D36") '0.06518'
$ws.Range("E36").Value = '  +6.67%  '

Set-TextValue $ws.Range("D37") '0.02254'
$ws.Range("E37").Value = '  -1.57%  '

Set-TextValue $ws.Range("D38") '1.194'
$ws.Range("E38").Value = '  +1.24%  '

Set-TextValue $ws.Range("D39") '0.5862'
$ws.Range("E39").Value = '  -2.16%  '

Set-TextValue $ws.Range("D40") '10.63'
$ws.Range("E40").Value = '  -1.49%  '

Set-TextValue $ws.Range("D41") '7.845'
$ws.Range("E41").Value = '  -2.78%  '

Set-TextValue $ws.Range("D42") '0.1827'
$ws.Range("E42").Value = '  -1.28%  '

Set-TextValue $ws.Range("D43") '2.460'
$ws.Range("E43").Value = '  +1.47%  '

Set-TextValue $ws.Range("D44") '1.276'
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D45") '0.07445'
$ws.Range("E45").Value = '  -2.37%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D46") '12.19'
$ws.Range("E46").Value = '  -2.48%  '

Set-TextValue $ws.Range("D47") '0.5501'
$ws.Range("E47").Value = '  -1.60%  '

Set-TextValue $ws.Range("D48") '1.963'
$ws.Range("E48").Value = '  +0.26%  '

Set-TextValue $ws.Range("D49") '116.15'
$ws.Range("E49").Value = '  -1.30%  '

$ws.Range("E50").Value = '  -1.21%  '

Set-TextValue $ws.Range("D51") '71.77'
$ws.Range("E51").Value = '  -0.95%  '
